$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "input" in column C for existing rows 4-7 (previously blank)
$ws.Range("C4").Value = "input"
$ws.Range("C5").Value = "input"
$ws.Range("C6").Value = "input"
$ws.Range("C7").Value = "input"

# Replace the long comment previously in C8 with "input"
$ws.Range("C8").Value = "input"

# Add new row 9 for historic_precipitation
$ws.Range("A9").Value = "historic_precipitation"
$ws.Range("B9").Value = "data/Dayun_weather_1-11-19_12-00_AM_1_Year_1583925430_v2.csv"
$ws.Range("C9").Value = 'Optional. In order to run the model with historical data, execute the program in the command line as "python daily_map -hp". The model will run with all the historical data contined in the file specified in Path. In order to change the range of the simulation, carefully edit the .csv file to the range required. It is crucial to maintain the format. When doing so, the water table  begins the simulation at the surface, i.e., in fully saturated conditions.'

# Update selection to C9 to match the final state
$ws.Range("C9").Select()
